# Build site at 2022-09-26 16:07:08 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (only B13/C13 populated with the docent name, no label in
# A13) is removed; everything below shifts up one row.
$ws.Rows("13").Delete()

# After the shift, re-point several label rows at their (new) neighbour's
# content, matching the edited spreadsheet.
$ws.Range("B10").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Range("C10").Value = "8188658 - Maria Auxiliadora Motta Barreto"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2012" must land as literal text (it already exists as a text
# shared string in B8/C8), so copy that cell instead of assigning .Value,
# which would otherwise auto-convert the look-alike date string to a
# serial date.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "8188658 - Maria Auxiliadora Motta Barreto"
$ws.Range("C18").Value = "8188658 - Maria Auxiliadora Motta Barreto"

$ws.Range("B19").Value = "utilizar-se-á provas dissertativas com estudo de caso, para levar os alunos à maior reflexão sobre a utilização dessa ciência para o futuro engenheiro em seu trabalho cotidiano em empresas."
$ws.Range("C19").Value = "utilizar-se-á provas dissertativas com estudo de caso, para levar os alunos à maior reflexão sobre a utilização dessa ciência para o futuro engenheiro em seu trabalho cotidiano em empresas."

$ws.Range("B20").Value = "(P1+ P2) : 2 = Média."
$ws.Range("C20").Value = "(P1+ P2) : 2 = Média."

$ws.Range("B21").Value = "Trabalho e nova avaliação para recuperação da nota necessária para aprovação."
$ws.Range("C21").Value = "Trabalho e nova avaliação para recuperação da nota necessária para aprovação."
